$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before assigning values so that
# numeric-looking strings (e.g. "1.00", "11.20") keep their exact text
# representation instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "42.297.92"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.250.06"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "247.62"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "74.45"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "41.92"
$ws.Range("E10").Value = "  +5.81%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  -3.71%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "7.18"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "14.51"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.856"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.254.06"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "42.160.19"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0982"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "71.97"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "2.26"
$ws.Range("E21").Value = "  +4.38%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "230.31"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "8.71"
$ws.Range("E23").Value = "  +36.38%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "11.20"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.63"
$ws.Range("E26").Value = "  -6.18%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.31"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "171.77"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.10"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "20.98"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0824"
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -4.88%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "30.36"
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.13"
$ws.Range("E35").Value = "  +7.40%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0307"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "13.44"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "2.19"
$ws.Range("E39").Value = "  -4.48%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.203"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "61.36"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "108.13"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "8.63"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.101"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "1.16"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "4.11"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  +0.35%  "
